# ExcelReporter: Align showing declared with showing detected licenses.
#
# For detected licenses only the license expression is shown, while for
# declared licenses the whole `ResolvedLicense` is shown. Align that by
# showing only the expression of declared licenses on the per-project sheet.

$wb = $excel.ActiveWorkbook

# The per-project sheet ("Gradle org.ossreviewtoolkit.gra") holds the
# "Declared Licenses" values as the whole ResolvedLicense(...) string.
# Replace them with just the license expression, like the Summary sheet
# already does.
$ws = $wb.Worksheets.Item("Gradle org.ossreviewtoolkit.gra")

$ws.Range("C12").Value = "EPL-1.0"
$ws.Range("C13").Value = "Apache-2.0"
$ws.Range("C14").Value = "Apache-2.0"
$ws.Range("C15").Value = "BSD-3-Clause"

$ws.Select()
$ws.Range("C15").Select()

$wb.Save()
